$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.841.84'
$ws.Range("E2").Value = '  +0.46%  '

$ws.Range("D3").Value = '2.462.20'
$ws.Range("E3").Value = '  +0.11%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '558.96'
$ws.Range("E5").Value = '  -0.63%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '161.80'
$ws.Range("E6").Value = '  -1.40%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.508'
$ws.Range("E8").Value = '  +0.41%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.151'
$ws.Range("E9").Value = '  -0.82%  '

$ws.Range("E10").Value = '  +0.68%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.331'
$ws.Range("E11").Value = '  -2.62%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.85'
$ws.Range("E12").Value = '  +0.30%  '

$ws.Range("D13").Value = '68.716.95'
$ws.Range("E13").Value = '  +0.38%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000168'
$ws.Range("E14").Value = '  -2.11%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.40'
$ws.Range("E15").Value = '  -0.73%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '10.63'
$ws.Range("E16").Value = '  -3.44%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '336.44'
$ws.Range("E17").Value = '  -2.06%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.92'
$ws.Range("E18").Value = '  -3.62%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.77'
$ws.Range("E19").Value = '  -1.27%  '

$ws.Range("B20").Value = 'SuiNetwork'
$ws.Range("C20").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.89'
$ws.Range("E20").Value = '  +0.83%  '

$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.00'
$ws.Range("E21").Value = '  -0.03%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '66.42'
$ws.Range("E22").Value = '  -2.40%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.64'
$ws.Range("E23").Value = '  -2.85%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.15'
$ws.Range("E24").Value = '  -1.00%  '

$ws.Range("D25").Value = '0.0₃0817'
$ws.Range("E25").Value = '  -2.93%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.20'
$ws.Range("E26").Value = '  -1.60%  '

$ws.Range("E27").Value = '  -0.02%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '429.66'
$ws.Range("E28").Value = '  -1.51%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.14'
$ws.Range("E29").Value = '  -3.60%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.61'
$ws.Range("E30").Value = '  -4.40%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '159.03'
$ws.Range("E31").Value = '  +1.21%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '19.02'
$ws.Range("E32").Value = '  -0.02%  '

$ws.Range("E33").Value = '  -0.03%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.108'
$ws.Range("E34").Value = '  -1.58%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '17.76'
$ws.Range("E35").Value = '  -0.92%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.300'
$ws.Range("E36").Value = '  -2.50%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.41'
$ws.Range("E37").Value = '  -2.05%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.46'
$ws.Range("E38").Value = '  -4.55%  '

$ws.Range("E39").Value = '  -3.18%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.05'
$ws.Range("E40").Value = '  -1.42%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.34'
$ws.Range("E41").Value = '  -1.04%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '129.52'
$ws.Range("E42").Value = '  -3.85%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0718'
$ws.Range("E43").Value = '  -0.15%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.482'
$ws.Range("E44").Value = '  -1.16%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.562'
$ws.Range("E45").Value = '  -0.06%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0911'
$ws.Range("E46").Value = '  -0.24%  '

$ws.Range("E47").Value = '  +0.18%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.38'
$ws.Range("E48").Value = '  -3.73%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.94'
$ws.Range("E49").Value = '  -8.42%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '16.78'
$ws.Range("E50").Value = '  -5.14%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.586'
$ws.Range("E51").Value = '  -3.47%  '
